$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.118.81'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '1.558.74'
$ws.Range('E3').Value = '  -0.26%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '209.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.489'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '21.96'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.70%  '
$ws.Range('E9').Value = '  -0.22%  '
$ws.Range('E10').Value = '  -1.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0869'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.56%  '
$ws.Range('D12').Value = '1.781.56'
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('D13').Value = '1.543.71'
$ws.Range('E13').Value = '  -1.26%  '
$ws.Range('E14').Value = '  +0.21%  '
$ws.Range('E15').Value = '  -0.78%  '
$ws.Range('D16').Value = '27.106.45'
$ws.Range('E16').Value = '  +0.50%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.65'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.54%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.43'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.82%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '216.12'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').Value = '0.0₃0700'
$ws.Range('E20').Value = '  -0.77%  '
$ws.Range('E21').Value = '  -0.12%  '
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.18'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.93'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.94'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.61'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.106'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.41%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '14.97'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.87%  '
$ws.Range('E30').Value = '  +1.81%  '
$ws.Range('E31').Value = '  -0.45%  '
$ws.Range('E32').Value = '  -0.25%  '
$ws.Range('E33').Value = '  +0.98%  '
$ws.Range('D34').Value = '1.433.05'
$ws.Range('E34').Value = '  +0.77%  '
$ws.Range('E35').Value = '  +3.61%  '
$ws.Range('E36').Value = '  -0.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.33'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.60%  '
$ws.Range('E38').Value = '  +0.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.529'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.91'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.39%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.804'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.49%  '
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.33'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.998'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.80%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.09'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.85%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.72'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.35%  '
$ws.Range('D47').Value = '1.695.11'
$ws.Range('E47').Value = '  -0.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '85.35'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0523'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.53%  '
$ws.Range('D50').Value = '0.0₇0992'
$ws.Range('E50').Value = '  -0.39%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0949'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.01%  '
